$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial value of 45186 (2023-09-17) for
# every data row (2 through 506). The commit updates that value to 45188
# (2023-09-19) for all of those rows, leaving everything else untouched.
$ws.Range("C2:C506").Value = 45188
